$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old used range entirely (A1:G4) so stale columns F/G are removed
$ws.Range("A1:G4").Clear()

# Headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 1057, 1057, 1057, 0.008758282661437989),
    @(1, 1040, 1040, 1040, 0.008761763572692871),
    @(2, 973, 973, 973, 0.01014242966969808),
    @(3, 1224, 1224, 1224, 0.01018432776133219),
    @(4, 883, 883, 883, 0.01228616237640381),
    @(5, 1040, 1040, 1040, 0.01369312604268392),
    @(6, 1053, 1053, 1053, 0.01368667284647624),
    @(7, 957, 957, 957, 0.01189967791239421),
    @(8, 886, 886, 886, 0.01372597217559814),
    @(9, 1049, 1049, 1049, 0.01384061177571615)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
